$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; existing data (A:E) shifts right to (B:F),
# carrying its formatting along with it.
$ws.Columns("A:A").Insert()

# Write the new header in B1 ("segments"), then restore the shared header style
# (same style used by C1:F1) by copying formats from C1, so Excel reuses the
# existing style index instead of allocating a brand-new one.
$ws.Range("B1").Value = "segments"
$ws.Range("C1").Copy() | Out-Null
$ws.Range("B1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Fill the new column A (rows 2-20) with the 0-based segment index.
for ($i = 0; $i -le 18; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $i
}

# Give the new index column (A2:A20) the same header-ish style that used to
# live on the (now shifted) name column, by copying format from B1.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("A2:A20").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# The name column (now B2:B20, previously A2:A20) keeps the old data but should
# no longer carry that bordered/bold style - clear it back to the default look.
$ws.Range("B2:B20").ClearFormats()
